$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before current row 2 (pushing existing data down)
$ws.Rows.Item(2).Resize(4).Insert()
$ws.Rows.Item(2).Resize(4).ClearFormats()

# New data to insert at the top (rows 2-5)
$topData = @(
    @(-0.0250454749912023, 0.0326812900602817, 0.0108428578823804),
    @(0.0192422550171613, -0.0277943685650825, -0.0045814891345798),
    @(-0.0242818929255008, 0.0114537235349416, 0.0073303831741213),
    @(0.0167987942695617, 0.1464549452066421, 0.5288565754890442)
)

for ($i = 0; $i -lt $topData.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $topData[$i][0]
    $ws.Cells.Item($r, 2).Value = $topData[$i][1]
    $ws.Cells.Item($r, 3).Value = $topData[$i][2]
}

# New data to append at the bottom (rows 26-31)
$bottomData = @(
    @(-0.1747074574232101, -0.5407684445381165, -0.012980886735022),
    @(-0.2935207486152649, -1.132696866989136, 0.0025961773935705),
    @(-0.1869247704744339, -0.6145304441452026, -0.2397646158933639),
    @(-0.1950187236070633, 0.0448985956609249, -0.3084869384765625),
    @(0.0444404482841491, 0.1437060534954071, -0.3220787048339844),
    @(0.2638937830924988, 0.3306308090686798, -0.3280346393585205)
)

for ($i = 0; $i -lt $bottomData.Length; $i++) {
    $r = 26 + $i
    $ws.Cells.Item($r, 1).Value = $bottomData[$i][0]
    $ws.Cells.Item($r, 2).Value = $bottomData[$i][1]
    $ws.Cells.Item($r, 3).Value = $bottomData[$i][2]
}
